$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells keep their original text data type,
# since several of the new values look like plain numbers (e.g. "1.01").
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "98.823.27"
$ws.Range("E2").Value = "  +0.81%  "
$ws.Range("D3").Value = "3.417.36"
$ws.Range("E3").Value = "  +3.75%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "259.07"
$ws.Range("E5").Value = "  +1.40%  "
$ws.Range("D6").Value = "671.02"
$ws.Range("E6").Value = "  +8.08%  "
$ws.Range("E7").Value = "  +10.28%  "
$ws.Range("E8").Value = "  +18.89%  "
$ws.Range("E9").Value = "  +22.64%  "
$ws.Range("E10").Value = "  -0.05%  "
$ws.Range("D11").Value = "3.415.48"
$ws.Range("E11").Value = "  +3.87%  "
$ws.Range("E12").Value = "  +11.13%  "
$ws.Range("D13").Value = "43.03"
$ws.Range("E13").Value = "  +13.63%  "
$ws.Range("E14").Value = "  +13.03%  "
$ws.Range("D15").Value = "6.14"
$ws.Range("E15").Value = "  +11.98%  "
$ws.Range("D16").Value = "98.304.13"
$ws.Range("E16").Value = "  +0.57%  "
$ws.Range("D17").Value = "4.057.91"
$ws.Range("E17").Value = "  +3.89%  "
$ws.Range("D18").Value = "8.18"
$ws.Range("E18").Value = "  +35.58%  "
$ws.Range("D19").Value = "3.418.74"
$ws.Range("E19").Value = "  +3.86%  "
$ws.Range("D20").Value = "17.35"
$ws.Range("E20").Value = "  +15.32%  "
$ws.Range("D21").Value = "537.12"
$ws.Range("E21").Value = "  +12.74%  "
$ws.Range("D22").Value = "3.59"
$ws.Range("E22").Value = "  +2.72%  "
$ws.Range("D23").Value = "10.64"
$ws.Range("E23").Value = "  +14.99%  "
$ws.Range("E24").Value = "  +8.58%  "
$ws.Range("D25").Value = "0.441"
$ws.Range("E25").Value = "  +53.09%  "
$ws.Range("E26").Value = "  +15.67%  "
$ws.Range("D27").Value = "103.08"
$ws.Range("E27").Value = "  +16.86%  "
$ws.Range("D28").Value = "12.90"
$ws.Range("E28").Value = "  +9.59%  "
$ws.Range("D29").Value = "3.597.21"
$ws.Range("E29").Value = "  +3.56%  "
$ws.Range("D30").Value = "0.152"
$ws.Range("E30").Value = "  +16.00%  "
$ws.Range("D31").Value = "11.68"
$ws.Range("E31").Value = "  +20.56%  "
$ws.Range("D32").Value = "0.198"
$ws.Range("E32").Value = "  +6.83%  "
$ws.Range("E33").Value = "  -0.11%  "
$ws.Range("D34").Value = "1.01"
$ws.Range("E34").Value = "  +0.55%  "
$ws.Range("D35").Value = "30.35"
$ws.Range("E35").Value = "  +10.84%  "
$ws.Range("E36").Value = "  +24.66%  "
$ws.Range("D37").Value = "2.21"
$ws.Range("E37").Value = "  +14.84%  "
$ws.Range("D38").Value = "7.92"
$ws.Range("E38").Value = "  +11.60%  "
$ws.Range("E39").Value = "  +10.43%  "
$ws.Range("D40").Value = "533.85"
$ws.Range("E40").Value = "  +8.76%  "
$ws.Range("D41").Value = "1.42"
$ws.Range("E41").Value = "  +15.57%  "
$ws.Range("E42").Value = "  -0.11%  "
$ws.Range("D43").Value = "0.0444"
$ws.Range("E43").Value = "  +38.24%  "
$ws.Range("D44").Value = "3.81"
$ws.Range("E44").Value = "  +4.00%  "
$ws.Range("E45").Value = "  +11.42%  "
$ws.Range("E46").Value = "  +8.92%  "
$ws.Range("E47").Value = "  +0.06%  "
$ws.Range("D48").Value = "8.19"
$ws.Range("E48").Value = "  +19.53%  "
$ws.Range("D49").Value = "2.12"
$ws.Range("E49").Value = "  +12.79%  "
$ws.Range("D50").Value = "5.28"
$ws.Range("E50").Value = "  +15.72%  "
$ws.Range("E51").Value = "  +17.45%  "
